# Actualizacion a 14 de Abril.
# Adds four new daily rows (2020-04-11 .. 2020-04-14, Excel serials 43931-43934)
# to both region sheets (Hoja1 = casos acumulados, Hoja2 = fallecidos acumulados),
# and updates each sheet's used range / view state to match.

$wb = $excel.ActiveWorkbook

# --- New data rows -----------------------------------------------------
# Columns: A fecha | B dia | C Arica y Parinacota | D Tarapaca | E Antofagasta
# F Atacama | G Coquimbo | H Valparaiso | I Metropolitana | J O'Higgins
# K Maule | L Nuble | M Biobio | N Araucania | O Los Rios | P Los Lagos
# Q Aysen | R Magallanes | S total

$sheet1Data = @{
    40 = @(43931,39,107,38,117,13,64,248,3448,45,134,606,490,739,130,364,7,377,6927)
    41 = @(43932,40,115,46,138,13,66,254,3599,46,138,613,500,775,135,372,7,396,7213)
    42 = @(43933,41,115,52,149,13,66,273,3803,48,141,618,512,795,138,380,7,415,7525)
    43 = @(43934,42,120,62,155,13,66,285,4086,53,142,622,528,816,148,385,7,429,7917)
}

$sheet2Data = @{
    40 = @(43931,39,0,0,1,0,0,2,32,0,3,6,2,17,2,5,0,3,73)
    41 = @(43932,40,1,0,1,0,0,2,35,0,3,7,2,17,3,5,0,4,80)
    42 = @(43933,41,1,0,1,0,0,2,36,0,3,7,2,17,3,5,0,5,82)
    43 = @(43934,42,1,0,1,0,0,2,40,0,4,8,2,20,3,5,0,6,92)
}

function Write-Rows($ws, $data) {
    foreach ($r in $data.Keys) {
        $vals = $data[$r]
        for ($c = 1; $c -le $vals.Length; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($c -eq 1) {
                # Column A keeps the same date format as the rest of the column
                $cell.NumberFormat = "DD/MM/YY"
            }
            $cell.Value = $vals[$c - 1]
        }
    }
}

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

Write-Rows $ws1 $sheet1Data
Write-Rows $ws2 $sheet2Data

# --- View / selection state ---------------------------------------------
# Hoja2 (not the active tab): scroll near the new rows and leave the
# selection on A41 (second area of the author's A52:AG54 + A41 selection).
[void]$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws2.Range("A41").Select()

# Hoja1 (the active tab): scroll down toward the new rows and select the
# block below the new data, matching the author's post-edit selection.
[void]$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws1.Range("A52:AG54").Select()
